$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row -> new text mapping (1-based row indices), per the target diff.
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "400"
$t.Cell(6, 1).Range.Text  = "0.00057"
$t.Cell(7, 1).Range.Text  = "0.00021"
$t.Cell(8, 1).Range.Text  = "0.00005"
$t.Cell(9, 1).Range.Text  = "0.00036"
$t.Cell(10, 1).Range.Text = "0.00040"
$t.Cell(11, 1).Range.Text = "0.00046"
$t.Cell(12, 1).Range.Text = "0.08526"

# These rows currently hold a tab-separated run list; collapse each to a
# single value (matching the earlier summary rows above).
$t.Cell(44, 1).Range.Text = "99.57"
$t.Cell(45, 1).Range.Text = "0.09"
$t.Cell(46, 1).Range.Text = "19"
